# The tracked change being replayed here lives entirely in PowerPoint's
# collaborative-editing change log (ppt/changesInfos/changesInfo1.xml):
# the top-level entry is flagged "undo custSel modSld" and the only
# shape-level entry touched is the title text box (shape id="186",
# "Architecture 3/3 - Car maintenance Expectation") on the slide with
# sldId="267" -- its chg flag grows from "mod" to "add del mod", i.e.
# the shape was deleted and that deletion was then undone in the same
# session. Because it's an undo, the slide's actual content/XML is left
# exactly as it was; only PowerPoint's internal revision-tracking
# bookkeeping (timestamps/version counters/action ids, none of which are
# reachable through the PowerPoint object model) records that the
# undo + re-selection happened.
#
# Reproduce the user-visible part of that action: locate the same slide
# and the same shape and (re-)select it, without altering its text,
# formatting or position.

$p = $ppt.ActivePresentation

$targetSlide = $null
foreach ($sl in $p.Slides) {
    if ($sl.SlideID -eq 267) {
        $targetSlide = $sl
        break
    }
}

if ($targetSlide -ne $null) {
    $targetShape = $null
    foreach ($shp in $targetSlide.Shapes) {
        if ($shp.Id -eq 186) {
            $targetShape = $shp
            break
        }
    }

    if ($targetShape -ne $null) {
        $targetShape.Select()
    }
}
